$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A62").NumberFormat = "@"
$ws.Range("D62").NumberFormat = "@"

$ws.Range("A62").Value = "2024-01-16"
$ws.Range("B62").Value = "09:44:17"
$ws.Range("C62").Value = "Tuesday"
$ws.Range("D62").Value = "02"
$ws.Range("E62").Value = 138806
$ws.Range("F62").Value = 139155
$ws.Range("G62").Value = 170793
$ws.Range("H62").Value = 148190
$ws.Range("I62").Value = -1
$ws.Range("J62").Value = 119006
$ws.Range("K62").Value = 220988
$ws.Range("L62").Value = 253956
$ws.Range("M62").Value = 184928
$ws.Range("N62").Value = 110378
$ws.Range("O62").Value = 41148
$ws.Range("P62").Value = 30880
$ws.Range("Q62").Value = 73188
$ws.Range("R62").Value = -1
$ws.Range("S62").Value = 42050
$ws.Range("T62").Value = -1
